$d = $word.ActiveDocument

# Locate the paragraph ending with the "如果一个类承担的职责过多..." sentence
# (last real content paragraph before the trailing empty/bookmark paragraph)
# by searching through the Paragraphs collection for the sentence.
$anchorIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "如果一个类承担的职责过多") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph for the Open/Closed Principle insertion."
}

$anchor = $d.Paragraphs.Item($anchorIndex)

# Insert a blank paragraph right after the anchor paragraph (mirrors the blank
# separator paragraph that precedes every new "N. <title>" section). This
# creates a brand-new empty paragraph right after the anchor.
$anchor.Range.InsertParagraphAfter()

# The newly created blank paragraph is now the paragraph right after the
# anchor. InsertXML *replaces* the contents of the range it targets, so we
# feed it the blank paragraph (re-created verbatim as the first <w:p>) plus
# the four new content paragraphs that should follow it - that way the
# blank separator paragraph survives and the new paragraphs land after it.
$target = $d.Paragraphs.Item($anchorIndex + 1)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">2. </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>开放封闭原则</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>类应该对扩展开放，对修改关闭。</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>扩展就是添加新功能的意思，因此该原则要求在添加新功能时不需要修改代码。</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>符合开闭原则最典型的设计模式是装饰者模式，它可以动态地将责任附加到对象上，而不用去修改类的代码。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($xml)
